# Fix mangled Vietnamese text (mojibake from a broken codepage round-trip)
# in the customer address column ("Địa chỉ") of Sheet1.
#
# "Ðu?ng" -> "Đường" and "Qu?n" -> "Quận" for the 3 sample customer rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "123 Đường ABC, Quận 1, TP.HCM"
$ws.Range("E3").Value = "456 Đường XYZ, Quận 2, TP.HCM"
$ws.Range("E4").Value = "789 Đường LMN, Quận 3, TP.HCM"
